# "Fix LLM Validation Result."
#
# The data row on "Main Sheet" (row 2) was written one column too far to the
# left starting at column D: the "Diagnosis Code" value was actually missing
# and every following field (Diagnosis Description, Provider Name, Provider
# Code, dates, amounts, checks, ...) had silently slid left by one column, so
# each value ended up misaligned under the wrong header. On top of that, the
# "Incur Date to" value had been duplicated into the "Provider Code" slot
# instead of getting its own cell.
#
# Fix: shift every populated cell from D2:R2 one column to the right (to
# E2:S2), which re-aligns all of the values with their correct headers and
# leaves D2 ("Diagnosis Code") empty, and insert a proper "Incur Date to"
# value in I2 (a copy of the "Incur Date from" value in H2, preserving its
# date number formatting) instead of the stray duplicate that used to live
# in G2.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Main Sheet")

# Shift the row, working from right to left so that a cell's original
# contents are always copied out before anything else overwrites it. Using
# Copy (rather than just copying Value2) also carries over the cell's style
# (e.g. the date number format on the "Incur Date" cells).
$ws1.Range("R2").Copy($ws1.Range("S2")) | Out-Null
$ws1.Range("Q2").Copy($ws1.Range("R2")) | Out-Null
$ws1.Range("P2").Copy($ws1.Range("Q2")) | Out-Null
$ws1.Range("O2").Copy($ws1.Range("P2")) | Out-Null
$ws1.Range("N2").Copy($ws1.Range("O2")) | Out-Null
$ws1.Range("M2").Copy($ws1.Range("N2")) | Out-Null
$ws1.Range("L2").Copy($ws1.Range("M2")) | Out-Null
$ws1.Range("K2").Copy($ws1.Range("L2")) | Out-Null
$ws1.Range("J2").Copy($ws1.Range("K2")) | Out-Null
$ws1.Range("I2").Copy($ws1.Range("J2")) | Out-Null
$ws1.Range("H2").Copy($ws1.Range("I2")) | Out-Null
$ws1.Range("F2").Copy($ws1.Range("G2")) | Out-Null
$ws1.Range("E2").Copy($ws1.Range("F2")) | Out-Null
$ws1.Range("D2").Copy($ws1.Range("E2")) | Out-Null
$ws1.Range("D2").Clear() | Out-Null

# Bring "Main Sheet" to the front (it was "Validation Summary" before) and
# leave the corrected row selected, matching the reviewed state of the file.
$ws1.Select() | Out-Null
$ws1.Range("E2:W2").Select() | Out-Null
